$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 93 (Vega Monumental Concepcion - Pina
# weekly price log). This shifts the existing row 93 -> 94 and row 94 -> 95,
# preserving their original values untouched.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with this week's data.
$ws.Range("A93").Value2 = 11
$ws.Range("B93").Value2 = "Vega Monumental Concepción"
$ws.Range("C93").Value2 = "Bíobío"
$ws.Range("D93").Value2 = 44448
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value2 = 8
$ws.Range("F93").Value2 = "Fruta"
$ws.Range("G93").Value2 = 100108
$ws.Range("H93").Value2 = "Tropicales y subtropicales"
$ws.Range("I93").Value2 = 100108005
$ws.Range("J93").Value2 = "Piña"
$ws.Range("K93").Value2 = "Caramelo"
$ws.Range("L93").Value2 = "Primera"
$ws.Range("M93").Value2 = 200
$ws.Range("N93").Value2 = 18500
$ws.Range("O93").Value2 = 19000
$ws.Range("P93").Value2 = 18750
$ws.Range("Q93").Value2 = "$/caja 12 unidades"
$ws.Range("R93").Value2 = "Ecuador"
$ws.Range("S93").Value2 = 1562
$ws.Range("T93").Value2 = 12
